$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 4 for columns D, M, N, O, P, R, S
$ws.Range("D2").Value = 44160
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 24500
$ws.Range("R2").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S2").Value = 1361

$ws.Range("D4").Value = 44174
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19500
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 1083
